$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential notice date in A9 from 2021-04-08 to 2021-04-09
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values
$ws.Range("D2").Value = 0.2459273407362129
$ws.Range("E2").Value = 0.008694556451612767

$ws.Range("D3").Value = 0.2462385614571746
$ws.Range("E3").Value = 0.008608321377331363

$ws.Range("D4").Value = 0.2580098933130877
$ws.Range("E4").Value = 0.009614015097564765

$ws.Range("D5").Value = 0.2498242044935249
$ws.Range("E5").Value = -0.0002599428125812375

$ws.Range("E6").Value = 0.006673500822807377
